$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("fishery_vuln_summary")
$ws.Activate()

# Headers (row 1) for the new summary block, columns O:S
$ws.Cells.Item(1, 15).Value = "Fishery"
$ws.Cells.Item(1, 16).Value = "Exposure"
$ws.Cells.Item(1, 17).Value = "Sensitivity"
$ws.Cells.Item(1, 18).Value = "Adaptive Capacity"
$ws.Cells.Item(1, 19).Value = "Vulnerability"

# Row 2 - CPS
$ws.Cells.Item(2, 15).Value = "CPS"
$ws.Cells.Item(2, 16).Value = 0.672619047619048
$ws.Cells.Item(2, 17).Value = 0.51587301587301604
$ws.Cells.Item(2, 18).Value = 0.52380952380952395
$ws.Cells.Item(2, 19).Value = 0.99260853199478505

# Row 3 - Dungeness crab (note capitalization fix vs. column B's "dungeness crab")
$ws.Cells.Item(3, 15).Value = "Dungeness crab"
$ws.Cells.Item(3, 16).Value = 0.62008928571428601
$ws.Cells.Item(3, 17).Value = 0.53622606449987398
$ws.Cells.Item(3, 18).Value = 0.46701388888888901
$ws.Cells.Item(3, 19).Value = 1.0145187856356299

# Row 4 - geoduck
$ws.Cells.Item(4, 15).Value = "geoduck"
$ws.Cells.Item(4, 16).Value = 0.73263888888888895
$ws.Cells.Item(4, 17).Value = 0.55208333333333304
$ws.Cells.Item(4, 18).Value = 0.421296296296296
$ws.Cells.Item(4, 19).Value = 1.09484434606767

# Row 5 - groundfish
$ws.Cells.Item(5, 15).Value = "groundfish"
$ws.Cells.Item(5, 16).Value = 0.65715811965811999
$ws.Cells.Item(5, 17).Value = 0.50278038847117801
$ws.Cells.Item(5, 18).Value = 0.43874643874643898
$ws.Cells.Item(5, 19).Value = 1.03197102180052

# Row 6 - hake
$ws.Cells.Item(6, 15).Value = "hake"
$ws.Cells.Item(6, 16).Value = 0.51388888888888895
$ws.Cells.Item(6, 17).Value = 0.407407407407407
$ws.Cells.Item(6, 18).Value = 0.52777777777777801
$ws.Cells.Item(6, 19).Value = 0.869403580132458

# Row 7 - HMS
$ws.Cells.Item(7, 15).Value = "HMS"
$ws.Cells.Item(7, 16).Value = 0.66866987179487203
$ws.Cells.Item(7, 17).Value = 0.51077915140415098
$ws.Cells.Item(7, 18).Value = 0.439727463312369
$ws.Cells.Item(7, 19).Value = 1.0371499655285901

# Row 8 - razor clam
$ws.Cells.Item(8, 15).Value = "razor clam"
$ws.Cells.Item(8, 16).Value = 0.52500000000000002
$ws.Cells.Item(8, 17).Value = 0.66666666666666696
$ws.Cells.Item(8, 18).Value = 0.43333333333333302
$ws.Cells.Item(8, 19).Value = 1.13508536198204

# Row 9 - salmon
$ws.Cells.Item(9, 15).Value = "salmon"
$ws.Cells.Item(9, 16).Value = 0.72720074268239399
$ws.Cells.Item(9, 17).Value = 0.53209012584012605
$ws.Cells.Item(9, 18).Value = 0.41815476190476197
$ws.Cells.Item(9, 19).Value = 1.1084252183673999

# Row 10 - scallops
$ws.Cells.Item(10, 15).Value = "scallops"
$ws.Cells.Item(10, 16).Value = 0.68452380952380998
$ws.Cells.Item(10, 17).Value = 0.52430555555555602
$ws.Cells.Item(10, 18).Value = 0.38888888888888901
$ws.Cells.Item(10, 19).Value = 1.0681751107167501

# Row 11 - sea cucumber
$ws.Cells.Item(11, 15).Value = "sea cucumber"
$ws.Cells.Item(11, 16).Value = 0.71347402597402598
$ws.Cells.Item(11, 17).Value = 0.54292929292929304
$ws.Cells.Item(11, 18).Value = 0.37373737373737398
$ws.Cells.Item(11, 19).Value = 1.1112084715113599

# Row 12 - sea urchin
$ws.Cells.Item(12, 15).Value = "sea urchin"
$ws.Cells.Item(12, 16).Value = 0.77440476190476204
$ws.Cells.Item(12, 17).Value = 0.55024509803921595
$ws.Cells.Item(12, 18).Value = 0.37581699346405201
$ws.Cells.Item(12, 19).Value = 1.1573814469053501

# Row 13 - shrimp
$ws.Cells.Item(13, 15).Value = "shrimp"
$ws.Cells.Item(13, 16).Value = 0.60551948051948101
$ws.Cells.Item(13, 17).Value = 0.52335858585858597
$ws.Cells.Item(13, 18).Value = 0.449494949494949
$ws.Cells.Item(13, 19).Value = 0.99727980518796699

# Row 14 - squid
$ws.Cells.Item(14, 15).Value = "squid"
$ws.Cells.Item(14, 16).Value = 0.63749999999999996
$ws.Cells.Item(14, 17).Value = 0.53888888888888897
$ws.Cells.Item(14, 18).Value = 0.55000000000000004
$ws.Cells.Item(14, 19).Value = 0.96075591244060798

# Apply the two-decimal number format to the new block, matching the style used for it
$ws.Range("O1:S14").NumberFormat = "0.00"

# Leave the new block selected, as the editor did before switching away
$ws.Range("O1:S14").Select()

# Add the new (blank) "Sheet2" tab at the end and make it the active sheet
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "Sheet2"
$newSheet.Activate()
